# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.340.87"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.99%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.573.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.66%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "507.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.23"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.90%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.579"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -5.69%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.577.71"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.55"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.36%  "
$ws.Range("E11").Value = "  -3.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.346"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.03%  "
$ws.Range("E13").Value = "  +0.84%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.029.39"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.09%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.387.58"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.60"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.69%  "
$ws.Range("E17").Value = "  -2.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.583.42"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "345.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.10"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.24%  "
$ws.Range("E23").Value = "  -0.19%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.419"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.91%  "
$ws.Range("E26").Value = "  -2.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.702.00"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.79%  "
$ws.Range("E28").Value = "  +0.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0₃0843"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.36"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.07%  "
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "19.36"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "152.68"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.88%  "
$ws.Range("E34").Value = "  -2.90%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.69"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.18"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.00%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.852"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.28%  "
$ws.Range("E39").Value = "  -2.36%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.842"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "36.12"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.31%  "
$ws.Range("E42").Value = "  -1.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "295.82"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -4.05%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.620"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.96%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0996"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.77%  "
$ws.Range("B46").Value = "FirstDigitalUSD"
$ws.Range("C46").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.998"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.18%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0556"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.66%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "19.66"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.84"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0232"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.58%  "
$ws.Range("E51").Value = "  -0.06%  "
